$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the "panel_query_time" timestamps on the "data" sheet (F2:F8) ---
$dataSheet.Range("F2").Value = "2021-10-05 14:22:09.602100"
$dataSheet.Range("F3").Value = "2021-10-05 14:22:09.602110"
$dataSheet.Range("F4").Value = "2021-10-05 14:22:09.602113"
$dataSheet.Range("F5").Value = "2021-10-05 14:22:09.602131"
$dataSheet.Range("F6").Value = "2021-10-05 14:22:09.602134"
$dataSheet.Range("F7").Value = "2021-10-05 14:22:09.602136"
$dataSheet.Range("F8").Value = "2021-10-05 14:22:09.602139"

# --- Add a new "metadata" sheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

# Header row B1:F1 - copy straight from the "data" sheet header so the
# (bold / bordered / centered) header style is reused exactly.
$dataSheet.Range("B1:F1").Copy($meta.Range("B1:F1"))

# G1 needs the same header style but isn't present on the data sheet -
# copy the style from B1 then overwrite the text.
$dataSheet.Range("B1").Copy($meta.Range("G1"))

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# A2 - copy style (+ the "0" index value) straight from the data sheet's A2.
$dataSheet.Range("A2").Copy($meta.Range("A2"))

# Data row values
$meta.Range("B2").Value = "Parathyroid Cancer"
$meta.Range("C2").Value = 86

# D2 must hold the literal text "1.4" (not the number 1.4) with NO special
# cell style. Stash the text in a scratch cell (stamping it with a text
# number-format so it is stored as a string), then paste-special just the
# VALUE into D2 so the destination's default (unstyled) format is kept.
$meta.Range("Z1").NumberFormat = "@"
$meta.Range("Z1").Value = "1.4"
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)  # xlPasteValues
$meta.Range("Z1").Clear()

$meta.Range("E2").Value = "2021-07-28T13:53:20.323052Z"
$meta.Range("F2").Value = "2021-10-05 14:22:09.598549"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/86/?format=json"

# Keep the original active sheet/selection on "data" (matches the source file).
$dataSheet.Activate() | Out-Null
$dataSheet.Range("A1").Select() | Out-Null
